# Generate Report for Handoff
# Replace the two e2e test file UUIDs with new ones, flip their status from
# "Handed back: in sync with en-US" to "Ready for handoff", and refresh the
# handoff/handback bookkeeping columns on the per-locale sheets to reflect
# that only the first file has produced a (shared) xliff so far.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "4930cd6e-3d3b-4858-a9c5-80fc33b1aac0"
$newUuid1 = "9a023f4c-dea3-4eab-91ea-79f72a209048"
$oldUuid2 = "c4afc7f6-4f30-4781-b948-86c6556fe580"
$newUuid2 = "ffffefd082a9-af3c-47af-88e1-af967a44b92d"

$oldStatus = "Handed back: in sync with en-US"
$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newUuid1.md"
$wsOverview.Range("A3").Value = "$newUuid2.md"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Range("G2").Value = "2016-08-20 13:03:55"
$wsOverview.Range("G3").Value = "2016-08-20 13:03:55"

# The hyperlinked cell's text and the hyperlink's own display text are two
# independent pieces of state here - update both so they stay in sync.
$wsOverview.Range("B2").Value = "e2e\$newUuid1.md"
$wsOverview.Range("B3").Value = "e2e\$newUuid2.md"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newUuid1.md"
    } elseif ($addr -eq '$B$3') {
        $hl.TextToDisplay = "e2e\$newUuid2.md"
    }
}

$wsOverview.Columns("E:F").ColumnWidth = 17.22

# ---------------------------------------------------------------------
# Helper to edit a per-locale sheet (zh-cn / de-de)
# ---------------------------------------------------------------------
function Update-LocaleSheet($ws, $xlfSuffix, $row2HandoffDate, $row3HandoffDate) {

    $ws.Range("A2").Value = "$newUuid1.md"
    $ws.Range("A3").Value = "$newUuid2.md"

    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    $xlfName = "$newUuid1.074d7c1959bb106be22360d7b6cb090df0a970a1.$xlfSuffix.xlf"

    # Row 2 - latest handoff/target info now points at the regenerated xlf.
    $ws.Range("G2").Value = $xlfName
    $ws.Range("H2").Value = $row2HandoffDate
    $ws.Range("I2").Value = ""
    $ws.Range("J2").Value = ""
    $ws.Range("K2").Value = "0001-01-01 00:00:00"

    # Row 3 - now a content duplicate of row 2, sharing the same xliff.
    # Leading apostrophe forces text (otherwise "True"/"False" literals get
    # auto-coerced to a boolean cell, unlike the rest of this text column).
    $ws.Range("F3").Value = "'True"
    $ws.Range("G3").Value = $xlfName
    $ws.Range("H3").Value = $row3HandoffDate
    $ws.Range("I3").Value = ""
    $ws.Range("J3").Value = ""
    $ws.Range("K3").Value = "0001-01-01 00:00:00"

    # Drop the now-stale "Latest Target File" hyperlinks on I2/I3, and reset
    # the cell style back to Normal (it carried the Hyperlink cell style).
    function Remove-HyperlinkAt($addr) {
        foreach ($hl in $ws.Hyperlinks) {
            if ($hl.Range.Address() -eq $addr) {
                $hl.Delete()
                return
            }
        }
    }
    Remove-HyperlinkAt '$I$2'
    Remove-HyperlinkAt '$I$3'
    $ws.Range("I2").Style = "Normal"
    $ws.Range("I3").Style = "Normal"

    # Update the still-present hyperlink display text for A2/A3.
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq '$A$2') {
            $hl.TextToDisplay = "$newUuid1.md"
        } elseif ($addr -eq '$A$3') {
            $hl.TextToDisplay = "$newUuid2.md"
        }
    }

    $ws.Columns("C:C").ColumnWidth = 17.22
    $ws.Columns("I:I").ColumnWidth = 18.65
    $ws.Columns("J:J").ColumnWidth = 21.71
}

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
Update-LocaleSheet $wsZh "zh-cn" "2016-08-20 13:03:51" "2016-08-20 13:03:51"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
Update-LocaleSheet $wsDe "de-de" "2016-08-20 13:03:55" "2016-08-20 13:03:55"
